# Add a new "Settings" worksheet after the existing "Sheet1", populate it
# with a small Key/Value table, and make it the active sheet (matching the
# commit "Settings Sheet added to .xlsx").

$wb = $excel.ActiveWorkbook

# Insert the new sheet after the last existing worksheet so it ends up
# positioned after Sheet1 (rather than Excel's default of inserting before
# the active sheet).
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "Settings"

$ws.Cells.Item(1, 1).Value = "Key"
$ws.Cells.Item(1, 2).Value = "Value"
$ws.Cells.Item(2, 1).Value = "League Champions Pot"
$ws.Cells.Item(2, 2).Value = 188.69

# Apply the default "Normal" cell style across the used range, as the
# source workbook does for every populated cell.
$ws.Range("A1:B2").Style = "Normal"

# Make the new Settings sheet the active tab, as in the target workbook.
$ws.Activate()
